$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A12").Style = "Normal"
$ws.Range("A12").Value = 3
$ws.Range("Z50").Formula = '=COUNTIF(W3:BG3,CONCATENATE("=",TEXT($A12,"d")))'
